# "EMail API added after testing"
# Update the existing customer-id value in A13 and append two newly
# generated customer ids in A14/A15, then move the active selection
# onto the newly added rows (A13:A15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "cus_KQ6aaMpTMkZP9V"
$ws.Range("A14").Value = "cus_KQ6ayxL6jsGJ8c"
$ws.Range("A15").Value = "cus_KQ6aP84xmpHq5P"

$ws.Range("A13:A15").Select()
